$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare formatting for the two brand-new rows (9 and 10) by copying the
# look of an existing data row (row 7) before we touch row 8's own content. ---
$ws.Range("A7:G7").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A10:G10").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 8: was the "Pause the script for 5 sec" case, now becomes the
# "Navigate to selenium download" case. (A8 keeps the same TC_007 text.) ---
$ws.Range("D8").Value = "Navigate"
$ws.Range("G8").Value = "https://www.selenium.dev/downloads/"
$ws.Range("C8").Value = "Navigate to selenium download"

# --- New row 9: download the selenium jar file ---
$ws.Range("A9").Value = "TC_008"
$ws.Range("B9").Value = "Multiple clicks on website"
$ws.Range("C9").Value = "Download jar file"
$ws.Range("D9").Value = "Click"
$ws.Range("E9").Value = "xpath"
$ws.Range("F9").Value = "//p[contains(text(),'Latest stable version')]//a[contains(text(),'3.141.59')]"

# --- New row 10: the old "pause" case, relocated here ---
$ws.Range("A10").Value = "TC_009"
$ws.Range("B10").Value = "Multiple clicks on website"
$ws.Range("C10").Value = "Pause"
$ws.Range("D10").Value = "Pause"
$ws.Range("G10").Value = 30

# --- Fix casing typo on the keyword used in row 2 (must be applied last so
# that the shared-string table keeps the same ordering as the source file) ---
$ws.Range("D2").Value = "ElementContainsText"

# --- Turn G8 into a real hyperlink, then restore its original cell style so
# the "Hyperlink" look-and-feel introduced by Excel is not kept on the cell ---
$ws.Hyperlinks.Add($ws.Range("G8"), "https://www.selenium.dev/downloads/")
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$wb.Styles.Item("Hyperlink").Delete()

# --- Widen columns C and G to fit the new, longer content (best-fit to the
# longest entries now present in each column) ---
$ws.Columns("C").ColumnWidth = 37
$ws.Columns("G").ColumnWidth = 45.25

# --- Update selection to match the author's final cursor position ---
$ws.Range("C9").Select()
